$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write the brand-new strings first, in the exact order they should be
# interned into the shared-string table (matches the authored file's table order). ---
$ws.Range("A3").Value = "Data driven testing"
$ws.Range("A1").Value = "Web Application/Topic"
$ws.Range("B3").Value = "Read data from excel sheet"
$ws.Range("B4").Value = "Write data in excel sheet"
$ws.Range("D4").Value = "In-Progress"
$ws.Range("E5").Value = "Read the values from website and put into the lists.`nNeed to write into the excel sheet."

# --- Step 2: write the remaining cells (these reuse already-interned strings). ---
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Test case"
$ws.Range("D1").Value = "Status"

$ws.Range("A2").Value = "Youtube"
$ws.Range("B2").Value = "Print all the video title and channel names"
$ws.Range("D2").Value = "Done"

$ws.Range("A4").Value = "Data driven testing"

$ws.Range("A5").Value = "Covid19India"
$ws.Range("B5").Value = "Get all states' corona count"
$ws.Range("C5").Value = "1.Go to https://www.covid19india.org/`n2.Get all states's Confirmed, Active, Recovered, Deceased count`n3.Print all the count in excel sheet"

$ws.Range("D3").Value = "Done"
$ws.Range("D5").Value = "In-Progress"

# The old row 3 used to carry the Covid description in column C; that data now
# lives in row 5, so clear the stale cell left behind in the new row 3.
$ws.Range("C3").Clear()

# --- Step 3: formatting for the wrapped description cells. ---
$ws.Range("C5").WrapText = $true
$ws.Range("C5").VerticalAlignment = -4160
$ws.Range("E5").WrapText = $true
$ws.Range("E5").VerticalAlignment = -4160

# --- Step 4: row heights. Rows 3 & 4 are short single-line rows, so drop any
# stale height inherited from the row they replaced; row 5 needs the taller,
# wrapped-text height (three lines @ 14.4pt). ---
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).RowHeight = 43.2

# --- Step 5: column E is now the wide "extra notes" column. ---
$ws.Columns.Item(5).ColumnWidth = 44.33

# --- Step 6: selection moves to A3. ---
$ws.Range("A3").Select()

# --- Step 7: data validation list now covers D2 through D19 (was D2:D20). ---
$ws.Range("D2:D20").Validation.Delete()
$ws.Range("D2:D19").Validation.Add(3, 1, 1, '"Done,In-Progress"')
